# Rename the workbook's only sheet from "Sheet1" to "Eagle"
# (input files updated to the Eagle pond per approval from parks)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Name = "Eagle"
